# Update ticket/sales count figures (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (1st tab / sheet1.xml) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 192
$ws1.Range("F5").Value = 1819
$ws1.Range("F6").Value = 467
$ws1.Range("F8").Value = 164
$ws1.Range("F9").Value = 2354
$ws1.Range("F10").Value = 132
$ws1.Range("F11").Value = 72
$ws1.Range("F12").Value = 155
$ws1.Range("F13").Value = 1428
$ws1.Range("F14").Value = 505
$ws1.Range("F15").Value = 34
$ws1.Range("F16").Value = 311
$ws1.Range("F17").Value = 221
$ws1.Range("F20").Value = 195
$ws1.Range("F21").Value = 210
$ws1.Range("F22").Value = 210
$ws1.Range("F24").Value = 92
$ws1.Range("F25").Value = 33
$ws1.Range("F26").Value = 1471
$ws1.Range("F28").Value = 370
$ws1.Range("F29").Value = 227
$ws1.Range("F32").Value = 370

# --- Sheet "全部类型" (4th tab / sheet4.xml) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 192
$ws4.Range("F5").Value = 1819
$ws4.Range("F7").Value = 467
$ws4.Range("F9").Value = 164
$ws4.Range("F10").Value = 2354
$ws4.Range("F11").Value = 132
$ws4.Range("F12").Value = 72
$ws4.Range("F13").Value = 155
$ws4.Range("F14").Value = 1428
$ws4.Range("F15").Value = 505
$ws4.Range("F16").Value = 34
$ws4.Range("F17").Value = 311
$ws4.Range("F18").Value = 221
$ws4.Range("F21").Value = 195
$ws4.Range("F22").Value = 210
$ws4.Range("F23").Value = 210
$ws4.Range("F25").Value = 92
$ws4.Range("F26").Value = 33
$ws4.Range("F27").Value = 1471
$ws4.Range("F29").Value = 370
$ws4.Range("F30").Value = 227
$ws4.Range("F33").Value = 370
